$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Values that Excel would otherwise auto-convert to a number (pure
    # integer/decimal strings) need a quote-prefix so they stay text -
    # exactly like the original file stores them (inline/shared string).
    if ($text -match '^[+-]?(\d+\.?\d*|\.\d+)$') {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '69.821.49'
$ws.Range("E2").Value = '  +2.83%  '
Set-TextValue $ws.Range("D3") '2.440.49'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.21%  '
Set-TextValue $ws.Range("D5") '564.63'
$ws.Range("E5").Value = '  +2.05%  '
Set-TextValue $ws.Range("D6") '166.10'
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.86%  '
Set-TextValue $ws.Range("D9") '0.170'
$ws.Range("E9").Value = '  +6.90%  '
Set-TextValue $ws.Range("D10") '2.444.46'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  -1.67%  '
Set-TextValue $ws.Range("D12") '0.334'
$ws.Range("E12").Value = '  +1.50%  '
Set-TextValue $ws.Range("D13") '4.70'
$ws.Range("E13").Value = '  -1.78%  '
Set-TextValue $ws.Range("D14") '0.0000179'
$ws.Range("E14").Value = '  +5.68%  '
Set-TextValue $ws.Range("D15") '69.860.02'
$ws.Range("E15").Value = '  +3.03%  '
Set-TextValue $ws.Range("D16") '2.896.82'
$ws.Range("E16").Value = '  -0.31%  '
Set-TextValue $ws.Range("D17") '23.99'
$ws.Range("E17").Value = '  +4.31%  '
Set-TextValue $ws.Range("D18") '2.446.94'
$ws.Range("E18").Value = '  +1.99%  '
Set-TextValue $ws.Range("D19") '10.78'
$ws.Range("E19").Value = '  +3.95%  '
Set-TextValue $ws.Range("D20") '340.26'
$ws.Range("E20").Value = '  +2.03%  '
Set-TextValue $ws.Range("D21") '7.11'
$ws.Range("E21").Value = '  +3.83%  '
Set-TextValue $ws.Range("D22") '3.88'
$ws.Range("E22").Value = '  +2.44%  '
Set-TextValue $ws.Range("D23") '1.99'
$ws.Range("E23").Value = '  +7.35%  '
Set-TextValue $ws.Range("D24") '0.999'
$ws.Range("E24").Value = '  -0.12%  '
Set-TextValue $ws.Range("D25") '66.07'
$ws.Range("E25").Value = '  -0.40%  '
Set-TextValue $ws.Range("D26") '3.80'
$ws.Range("E26").Value = '  +4.71%  '
Set-TextValue $ws.Range("D27") '2.571.76'
$ws.Range("E27").Value = '  +0.58%  '
Set-TextValue $ws.Range("D28") '8.45'
$ws.Range("E28").Value = '  +4.31%  '
Set-TextValue $ws.Range("D29") '0.996'
$ws.Range("E29").Value = '  -0.26%  '
Set-TextValue $ws.Range("D30") '0.0₃0850'
$ws.Range("E30").Value = '  +5.46%  '
Set-TextValue $ws.Range("D31") '7.34'
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("E32").Value = '  +9.35%  '
Set-TextValue $ws.Range("D33") '452.47'
$ws.Range("E33").Value = '  +8.78%  '
$ws.Range("E34").Value = '  +0.44%  '
Set-TextValue $ws.Range("D35") '1.61'
$ws.Range("E35").Value = '  +1.01%  '
Set-TextValue $ws.Range("D36") '159.45'
$ws.Range("E36").Value = '  -0.56%  '
Set-TextValue $ws.Range("D37") '19.08'
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("E38").Value = '  +5.04%  '
$ws.Range("E39").Value = '  +0.03%  '
Set-TextValue $ws.Range("D40") '18.18'
$ws.Range("E40").Value = '  +1.97%  '
Set-TextValue $ws.Range("D41") '0.305'
$ws.Range("E41").Value = '  +3.49%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D42") '38.02'
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D43") '1.52'
$ws.Range("E43").Value = '  +3.88%  '
Set-TextValue $ws.Range("D44") '4.38'
$ws.Range("E44").Value = '  +2.74%  '
Set-TextValue $ws.Range("D45") '1.08'
$ws.Range("E45").Value = '  +0.92%  '
Set-TextValue $ws.Range("D46") '2.10'
$ws.Range("E46").Value = '  +5.18%  '
Set-TextValue $ws.Range("D47") '134.53'
$ws.Range("E47").Value = '  +3.83%  '
Set-TextValue $ws.Range("D48") '3.38'
$ws.Range("E48").Value = '  +1.69%  '
Set-TextValue $ws.Range("D49") '0.0727'
$ws.Range("E49").Value = '  +3.06%  '
Set-TextValue $ws.Range("D50") '0.488'
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("E51").Value = '  +1.40%  '